# Auto-generated Excel COM-interop script to update crypto price/volume data
# Applies the GitHub Actions data refresh for cryptos.xlsx (Thu Jun 6 22:35:26 UTC 2024)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($sheet, $ref, $val) {
    # Force the cell to be written as literal text, matching the workbook's
    # existing inline-string cell type (so numeric-looking values like "170.52"
    # are not silently coerced into the Number type), then restore the cell's
    # original (default) style so no stray formatting is introduced.
    $range = $sheet.Range($ref)
    $range.NumberFormat = "@"
    $range.Value = $val
    $range.Style = "Normal"
}

Set-TextCell $ws 'D2' '70.729.71'
Set-TextCell $ws 'E2' '  -0.56%  '
Set-TextCell $ws 'D3' '3.806.56'
Set-TextCell $ws 'E3' '  -0.99%  '
Set-TextCell $ws 'E4' '  -0.09%  '
Set-TextCell $ws 'D5' '707.38'
Set-TextCell $ws 'D6' '170.52'
Set-TextCell $ws 'E6' '  -1.74%  '
Set-TextCell $ws 'D7' '3.805.93'
Set-TextCell $ws 'E7' '  -0.97%  '
Set-TextCell $ws 'E8' '  +0.04%  '
Set-TextCell $ws 'E9' '  -0.91%  '
Set-TextCell $ws 'E10' '  -1.36%  '
Set-TextCell $ws 'E11' '  +1.89%  '
Set-TextCell $ws 'D12' '0.457'
Set-TextCell $ws 'E12' '  -1.01%  '
Set-TextCell $ws 'D13' '0.0000253'
Set-TextCell $ws 'E13' '  -1.62%  '
Set-TextCell $ws 'D14' '36.08'
Set-TextCell $ws 'E14' '  -0.96%  '
Set-TextCell $ws 'D15' '4.448.51'
Set-TextCell $ws 'E15' '  -0.98%  '
Set-TextCell $ws 'D16' '3.817.66'
Set-TextCell $ws 'E16' '  -0.99%  '
Set-TextCell $ws 'D17' '70.762.45'
Set-TextCell $ws 'E17' '  -0.59%  '
Set-TextCell $ws 'E18' '  +0.14%  '
Set-TextCell $ws 'D19' '7.13'
Set-TextCell $ws 'E19' '  -1.38%  '
Set-TextCell $ws 'D20' '17.39'
Set-TextCell $ws 'E20' '  -1.86%  '
Set-TextCell $ws 'D21' '495.16'
Set-TextCell $ws 'E21' '  +0.39%  '
Set-TextCell $ws 'E22' '  -4.77%  '
Set-TextCell $ws 'E23' '  +0.98%  '
Set-TextCell $ws 'D24' '84.37'
Set-TextCell $ws 'E24' '  -0.73%  '
Set-TextCell $ws 'E25' '  -0.85%  '
Set-TextCell $ws 'E26' '  -1.72%  '
Set-TextCell $ws 'D27' '10.41'
Set-TextCell $ws 'E27' '  -1.70%  '
Set-TextCell $ws 'D28' '3.957.09'
Set-TextCell $ws 'E28' '  -1.19%  '
Set-TextCell $ws 'E29' '  +0.10%  '
Set-TextCell $ws 'E30' '  -4.82%  '
Set-TextCell $ws 'D31' '3.08'
Set-TextCell $ws 'E31' '  -2.99%  '
Set-TextCell $ws 'D32' '2.25'
Set-TextCell $ws 'E32' '  -1.63%  '
Set-TextCell $ws 'D33' '7.33'
Set-TextCell $ws 'E33' '  -3.91%  '
Set-TextCell $ws 'D34' '29.10'
Set-TextCell $ws 'E34' '  -1.85%  '
Set-TextCell $ws 'E35' '  -3.55%  '
Set-TextCell $ws 'B36' 'Aptos'
Set-TextCell $ws 'C36' 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextCell $ws 'D36' '9.13'
Set-TextCell $ws 'E36' '  -1.68%  '
Set-TextCell $ws 'B37' 'RenzoRestakedETH'
Set-TextCell $ws 'C37' 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
Set-TextCell $ws 'D37' '3.776.61'
Set-TextCell $ws 'E37' '  -0.53%  '
Set-TextCell $ws 'B38' 'Binance-PegBSC-USD'
Set-TextCell $ws 'C38' 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextCell $ws 'D38' '0.999'
Set-TextCell $ws 'E38' '  -0.15%  '
Set-TextCell $ws 'D39' '0.102'
Set-TextCell $ws 'E39' '  -3.04%  '
Set-TextCell $ws 'E40' '  +1.25%  '
Set-TextCell $ws 'D41' '2.31'
Set-TextCell $ws 'E41' '  -2.94%  '
Set-TextCell $ws 'D42' '5.94'
Set-TextCell $ws 'E42' '  -1.71%  '
Set-TextCell $ws 'E43' '  -3.76%  '
Set-TextCell $ws 'B45' 'FLOKI'
Set-TextCell $ws 'C45' 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
Set-TextCell $ws 'D45' '0.000324'
Set-TextCell $ws 'E45' '  +6.20%  '
Set-TextCell $ws 'B46' 'FirstDigitalUSD'
Set-TextCell $ws 'C46' 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextCell $ws 'D46' '1.00'
Set-TextCell $ws 'E46' '  +0.03%  '
Set-TextCell $ws 'D47' '165.01'
Set-TextCell $ws 'E47' '  +1.09%  '
Set-TextCell $ws 'D48' '425.03'
Set-TextCell $ws 'E48' '  +1.62%  '
Set-TextCell $ws 'D49' '48.75'
Set-TextCell $ws 'E49' '  +0.25%  '
Set-TextCell $ws 'E50' '  -0.46%  '
Set-TextCell $ws 'B51' 'TheGraph'
Set-TextCell $ws 'C51' 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextCell $ws 'D51' '0.294'
Set-TextCell $ws 'E51' '  -2.71%  '
